$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = 0.79018535775080501
$ws.Range("S1").Value = 0.69576715418328594
$ws.Range("D2").Value = 0.7031089418635007
$ws.Range("B3").Value = 0.88787566080525582
$ws.Range("AH3").Value = 0.84995011930361619
$ws.Range("AI4").Value = 0.9765819274626697
$ws.Range("C5").Value = 0.98449714208956496
$ws.Range("K5").Value = 0.95513095931378633
$ws.Range("W5").Value = 0.81100850698411664
$ws.Range("BK5").Value = 0.85526080320548659
$ws.Range("BM5").Value = 0.75300269193620206
$ws.Range("AB6").Value = 0.95018849545554251
$ws.Range("J8").Value = 0.89602084631394652
$ws.Range("AA8").Value = 0.73600656612360571
$ws.Range("J9").Value = 0.80626688197643537
$ws.Range("AB9").Value = 0.79076226732455523
$ws.Range("AD10").Value = 0.66511208963471757
$ws.Range("AO10").Value = 0.88937600476641099
$ws.Range("AW10").Value = 0.99899785976173505
$ws.Range("L11").Value = 0.87267717828443714
$ws.Range("Y12").Value = 0.72109625596507732
$ws.Range("AV13").Value = 0.93587917286807554
$ws.Range("H14").Value = 0.76019292256769089
$ws.Range("BK14").Value = 0.69054615329796176
$ws.Range("M15").Value = 0.86677883135664335
$ws.Range("BO15").Value = 0.83480588544593382
$ws.Range("D16").Value = 0.64758064153153905
$ws.Range("Q16").Value = 0.61936001823260201
$ws.Range("AE16").Value = 0.80325076403895723
$ws.Range("BL16").Value = 0.93150937489426089
$ws.Range("S17").Value = 0.756576339529201
$ws.Range("AS17").Value = 0.99917447557990213
$ws.Range("BH17").Value = 0.71565319465332133
$ws.Range("P18").Value = 0.90384115330057657
$ws.Range("W18").Value = 0.99071476038459161
$ws.Range("BE18").Value = 0.6623749018095999
$ws.Range("BC19").Value = 0.85288851592153669
$ws.Range("AC20").Value = 0.97640950058664844
$ws.Range("BN20").Value = 0.72940181077259358
$ws.Range("S21").Value = 0.97604307442460558
$ws.Range("AP21").Value = 0.956596073124097
$ws.Range("BN21").Value = 0.74044774567445004
$ws.Range("J22").Value = 0.89551807731388067
$ws.Range("Y24").Value = 0.83239868529799588
$ws.Range("AN24").Value = 0.64177570949874752
$ws.Range("E25").Value = 0.63627467639839552
$ws.Range("N25").Value = 0.94000264106789566
$ws.Range("BC25").Value = 0.6260054609077913
$ws.Range("BN25").Value = 0.64837670085287491
$ws.Range("Q26").Value = 0.76403612539546528
$ws.Range("S26").Value = 0.91887093951296794
$ws.Range("AA26").Value = 0.90384952706383559
$ws.Range("BG26").Value = 0.88329002693672387
$ws.Range("G28").Value = 0.94755420208512386
$ws.Range("AC28").Value = 0.83660340813047196
$ws.Range("AS29").Value = 0.79055575039330861
$ws.Range("BN29").Value = 0.8096935007882583
$ws.Range("Y30").Value = 0.79753030593918162
$ws.Range("AB30").Value = 0.81947018426648743
$ws.Range("AC30").Value = 0.86951798018655779
$ws.Range("AY30").Value = 0.95417480482254102
$ws.Range("BO31").Value = 0.92123005172910599
$ws.Range("F32").Value = 0.74544319674525106
$ws.Range("AP32").Value = 0.84941820417649283
$ws.Range("AY32").Value = 0.81722721437534629
$ws.Range("BO32").Value = 0.81858550890327386
$ws.Range("L33").Value = 0.92741919146586116
$ws.Range("AG34").Value = 0.87263251927219199
$ws.Range("AI34").Value = 0.75057818544240384
$ws.Range("I35").Value = 0.57331985238856387
$ws.Range("R35").Value = 0.86169397716978047
$ws.Range("AN35").Value = 0.77033594566730446
$ws.Range("AS35").Value = 0.96675497135287325
$ws.Range("BI35").Value = 0.85322746813735062
$ws.Range("AT36").Value = 0.88280646710731148
$ws.Range("AV36").Value = 0.82810492877370079
$ws.Range("AX36").Value = 0.91821054962143256
$ws.Range("BB37").Value = 0.88523117772581761
$ws.Range("BN37").Value = 0.9031959580633222
$ws.Range("AJ38").Value = 0.77007352570502285
$ws.Range("AK38").Value = 0.88683404505860641
$ws.Range("AS38").Value = 0.8878446392593462
$ws.Range("D39").Value = 0.81693578906598452
$ws.Range("E39").Value = 0.92650819652594851
$ws.Range("G40").Value = 0.82477703934416624
$ws.Range("Y40").Value = 0.75684374976597313
$ws.Range("AP40").Value = 0.94998423776447494
$ws.Range("Z41").Value = 0.65876693632065453
$ws.Range("AP41").Value = 0.95539783835504055
$ws.Range("BJ42").Value = 0.55313963743968619
$ws.Range("H43").Value = 0.70665628648062784
$ws.Range("BB43").Value = 0.79347437694639145
$ws.Range("BI43").Value = 0.70827170310600929
$ws.Range("A44").Value = 0.86604774456873379
$ws.Range("AI44").Value = 0.95518873983132346
$ws.Range("AR45").Value = 0.92595582390999875
$ws.Range("AU45").Value = 0.90801173144620162
$ws.Range("AI46").Value = 0.84433162863099764
$ws.Range("AZ47").Value = 0.78549630371845347
$ws.Range("BL47").Value = 0.88677509468632709
$ws.Range("P48").Value = 0.85962838291168964
$ws.Range("C49").Value = 0.88756943839004632
$ws.Range("L49").Value = 0.98076436028305713
$ws.Range("AH49").Value = 0.83689941199121631
$ws.Range("AU49").Value = 0.94329904372179141
$ws.Range("BG50").Value = 0.92561176583749494
$ws.Range("BE52").Value = 0.75828482042954226
$ws.Range("AA53").Value = 0.67053213841378145
$ws.Range("AM53").Value = 0.94809832096793945
$ws.Range("BC54").Value = 0.86332551907873656
$ws.Range("R55").Value = 0.92738685193224946
$ws.Range("U55").Value = 0.92143504173130286
$ws.Range("AJ55").Value = 0.96034134877972621
$ws.Range("A56").Value = 0.66179256081261295
$ws.Range("BG56").Value = 0.83243223527337618
$ws.Range("M57").Value = 0.93266157258777183
$ws.Range("AS57").Value = 0.75762307662870421
$ws.Range("BD57").Value = 0.97662872435796177
$ws.Range("C58").Value = 0.89284263218044668
$ws.Range("BD58").Value = 0.83448622805658856
$ws.Range("B60").Value = 0.95475615606438435
$ws.Range("AY61").Value = 0.69813040265485948
$ws.Range("BA61").Value = 0.80158359597109907
$ws.Range("BJ61").Value = 0.62129081195453684
$ws.Range("BL62").Value = 0.99842875449506108
$ws.Range("BG63").Value = 0.87482823961610501
$ws.Range("BJ63").Value = 0.79710547056723768
$ws.Range("AP64").Value = 0.89739070863922188
$ws.Range("V65").Value = 0.76391190988117863
$ws.Range("AA65").Value = 0.929645074848652
$ws.Range("AV66").Value = 0.65875399737478291
$ws.Range("H67").Value = 0.60042104541827612
$ws.Range("BP67").Value = 0.9181476582071556
$ws.Range("R68").Value = 0.70531692991587058
$ws.Range("BB68").Value = 0.75661221783515409
$ws.Range("BM68").Value = 0.99203902496058682
